$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @(1362.45, 1342.7, 1351.55, 1351.9, 92, 1358.85)
    3  = @(1121.85, 1091.15, 1095.75, 1097.8, 117, 1099)
    4  = @(47780, 47424.1, 47685.65, 47691.25, 20, 47703.65)
    5  = @(598.7, 587.45, 595.45, 595.9, 184, 593.85)
    6  = @(946.05, 929.55, 936.35, 935.6, 113, 931.75)
    7  = @(533.2, 522.35, 528.75, 528.55, 133, 526.2)
    8  = @(1097.8, 1087, 1094.5, 1094.8, 178, 1091.45)
    9  = @(852.75, 830.15, 843.15, 842.65, 86, 831)
    10 = @(22544.95, 22466, 22500.95, 22497.95, 37, 22517.45)
    11 = @(3037, 2995.1, 3022.1, 3026.45, 70, 3002.45)
    12 = @(779.8, 772.85, 774.95, 775.2, 244, 778.35)
    13 = @(1216.8, 1196.35, 1209.65, 1207.95, 20, 1202.95)
    14 = @(996.4, 985.4, 993.6, 992.5, 129, 996)
    15 = @(155.85, 151.8, 154.05, 154.05, 1317, 153.45)
    16 = @(4135.9, 4104, 4108.6, 4108.15, 14, 4133.45)
    17 = @(3773.95, 3735.9, 3749, 3748.6, 7, 3755.6)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
